# Auto-generated edit script applying numeric updates described by the diff
# against Sheets/Pandaemonium_Profits.xlsx (sheet tabs: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# Row 18 (ALC) - diff hunk @ -1529
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 387.41177
$ws.Range("I18").Value = 349.125
$ws.Range("K18").Value = 349.125
$ws.Range("M18").Value = -65.125

# Row 131 (ALC) - diff hunk @ -7174
$ws.Range("H131").Value = 3983.8572
$ws.Range("I131").Value = 610
$ws.Range("J131").Value = 6060.077
$ws.Range("K131").Value = 1830
$ws.Range("L131").Value = 18180.231
$ws.Range("M131").Value = 3210
$ws.Range("N131").Value = -28260.231

# Row 4 (ARM) - diff hunk @ -7926
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 100
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -332

# Row 106 (ARM) - diff hunk @ -12927
$ws.Range("H106").Value = 50370
$ws.Range("J106").Value = 50370
$ws.Range("L106").Value = 50370
$ws.Range("N106").Value = -52894

# Row 128 (ARM) - diff hunk @ -13993
$ws.Range("H128").Value = 74566.664
$ws.Range("J128").Value = 74566.664
$ws.Range("L128").Value = 74566.664
$ws.Range("N128").Value = -84526.664

# Row 132 (ARM) - diff hunk @ -14186
$ws.Range("H132").Value = 2028.9744
$ws.Range("I132").Value = 1821.6666
$ws.Range("K132").Value = 5464.9998
$ws.Range("M132").Value = -2934.9998

# Row 134 (BSM) - diff hunk @ -21193
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2677.08
$ws.Range("I134").Value = 2686.4
$ws.Range("K134").Value = 8059.200000000001
$ws.Range("M134").Value = -5524.200000000001

# Row 22 (CRP) - diff hunk @ -22665
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 241.86667
$ws.Range("I22").Value = 276.375
$ws.Range("J22").Value = 202.42857
$ws.Range("K22").Value = 276.375
$ws.Range("L22").Value = 202.42857
$ws.Range("M22").Value = 73.625
$ws.Range("N22").Value = -902.42857

# Row 58 (CRP) - diff hunk @ -24438
$ws.Range("H58").Value = 2335651.5
$ws.Range("I58").Value = 3368766.8
$ws.Range("K58").Value = 3368766.8
$ws.Range("M58").Value = -3368563.8

# Row 134 (CRP) - diff hunk @ -28189
$ws.Range("H134").Value = 2379.4062
$ws.Range("I134").Value = 2115.1482
$ws.Range("J134").Value = 3806.4
$ws.Range("K134").Value = 6345.444600000001
$ws.Range("L134").Value = 11419.2
$ws.Range("M134").Value = -3810.444600000001
$ws.Range("N134").Value = -16489.2

# Row 136 (CRP) - diff hunk @ -28290
$ws.Range("H136").Value = 2335651.5
$ws.Range("I136").Value = 3368766.8
$ws.Range("K136").Value = 10106300.4
$ws.Range("M136").Value = -10103750.4

# Row 138 (CRP) - diff hunk @ -28391
$ws.Range("H138").Value = 35880
$ws.Range("J138").Value = 35880
$ws.Range("L138").Value = 35880
$ws.Range("N138").Value = -46160

# Row 139 (CRP) - diff hunk @ -28437
$ws.Range("H139").Value = 47500
$ws.Range("J139").Value = 47500
$ws.Range("L139").Value = 47500
$ws.Range("N139").Value = -57780

# Row 2 (CUL) - diff hunk @ -28681
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 23.838709
$ws.Range("I2").Value = 59.5
$ws.Range("J2").Value = 21.37931
$ws.Range("K2").Value = 357
$ws.Range("L2").Value = 128.27586
$ws.Range("M2").Value = -244
$ws.Range("N2").Value = -354.27586

# Row 9 (CUL) - diff hunk @ -29042
$ws.Range("H9").Value = 27017.277
$ws.Range("J9").Value = 27017.277
$ws.Range("L9").Value = 81051.83099999999
$ws.Range("N9").Value = -81499.83099999999

# Row 39 (CUL) - diff hunk @ -30566
$ws.Range("H39").Value = 3866.641
$ws.Range("J39").Value = 4072.1943
$ws.Range("L39").Value = 12216.5829
$ws.Range("N39").Value = -12804.5829

# Row 109 (CUL) - diff hunk @ -34098
$ws.Range("H109").Value = 2216
$ws.Range("I109").Value = 1660
$ws.Range("J109").Value = 2772
$ws.Range("K109").Value = 4980
$ws.Range("L109").Value = 8316
$ws.Range("M109").Value = -3940
$ws.Range("N109").Value = -10396

# Row 131 (CUL) - diff hunk @ -35224
$ws.Range("H131").Value = 1484.4
$ws.Range("J131").Value = 1110.8
$ws.Range("L131").Value = 3332.4
$ws.Range("N131").Value = -13412.4

# Row 64 (GSM) - diff hunk @ -38937
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 37771
$ws.Range("J64").Value = 37771
$ws.Range("L64").Value = 37771
$ws.Range("N64").Value = -38267

# Row 67 (GSM) - diff hunk @ -39084
$ws.Range("H67").Value = 37771
$ws.Range("J67").Value = 37771
$ws.Range("L67").Value = 37771
$ws.Range("N67").Value = -39487

# Row 121 (GSM) - diff hunk @ -41709
$ws.Range("H121").Value = 18587.5
$ws.Range("J121").Value = 18587.5
$ws.Range("L121").Value = 18587.5
$ws.Range("N121").Value = -22081.5

# Row 124 (GSM) - diff hunk @ -41853
$ws.Range("H124").Value = 79800
$ws.Range("J124").Value = 79800
$ws.Range("L124").Value = 79800
$ws.Range("N124").Value = -89620

# Row 131 (GSM) - diff hunk @ -42196
$ws.Range("H131").Value = 44000
$ws.Range("J131").Value = 44000
$ws.Range("L131").Value = 44000
$ws.Range("N131").Value = -54080

# Row 22 (LTW) - diff hunk @ -43812
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 951
$ws.Range("J22").Value = 1402
$ws.Range("L22").Value = 1402
$ws.Range("N22").Value = -1992

# Row 27 (LTW) - diff hunk @ -44060
$ws.Range("H27").Value = 951
$ws.Range("J27").Value = 1402
$ws.Range("L27").Value = 1402
$ws.Range("N27").Value = -1616

# Row 46 (LTW) - diff hunk @ -44994
$ws.Range("H46").Value = 684.9259
$ws.Range("I46").Value = 903.5
$ws.Range("J46").Value = 622.4761999999999
$ws.Range("K46").Value = 903.5
$ws.Range("L46").Value = 622.4761999999999
$ws.Range("M46").Value = -715.5
$ws.Range("N46").Value = -998.4761999999999

# Row 100 (LTW) - diff hunk @ -47652
$ws.Range("H100").Value = 4862.5
$ws.Range("I100").Value = 4180
$ws.Range("K100").Value = 4180
$ws.Range("M100").Value = -3639

# Row 128 (LTW) - diff hunk @ -49021
$ws.Range("H128").Value = 70200
$ws.Range("J128").Value = 70200
$ws.Range("L128").Value = 70200
$ws.Range("N128").Value = -80160

# Row 132 (LTW) - diff hunk @ -49217
$ws.Range("H132").Value = 3744.2083
$ws.Range("I132").Value = 3226
$ws.Range("J132").Value = 5298.8335
$ws.Range("K132").Value = 9678
$ws.Range("L132").Value = 15896.5005
$ws.Range("M132").Value = -7148
$ws.Range("N132").Value = -20956.5005

# Row 10 (WVR) - diff hunk @ -50211
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 6000
$ws.Range("J10").Value = 6000
$ws.Range("L10").Value = 6000
$ws.Range("N10").Value = -6338

# Row 13 (WVR) - diff hunk @ -50361
$ws.Range("H13").Value = 3624
$ws.Range("I13").Value = 3624
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 3624
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -3484
$ws.Range("N13").ClearContents()

# Row 47 (WVR) - diff hunk @ -52027
$ws.Range("H47").Value = 175000
$ws.Range("J47").Value = 175000
$ws.Range("L47").Value = 175000
$ws.Range("N47").Value = -176144

# Row 131 (WVR) - diff hunk @ -56140
$ws.Range("H131").Value = 47838.332
$ws.Range("J131").Value = 47838.332
$ws.Range("L131").Value = 47838.332
$ws.Range("N131").Value = -57918.332

